$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new data row for 2020-09-15 (date 44089)
$ws.Range("A24").Value = 44089
$ws.Range("A24").NumberFormat = "d-mmm"

$ws.Range("B24").Value = 0.6875
$ws.Range("B24").NumberFormat = "h:mm"

$ws.Range("C24").Value = 0.89583333333333337
$ws.Range("C24").NumberFormat = "h:mm"

# Update the sheet view to reflect current selection/scroll position
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("E13").Select()
